$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $val
    $range.Style = "Normal"
}

Set-TextValue "D2" "66.224.12"
Set-TextValue "E2" "  -4.47%  "
Set-TextValue "D3" "3.343.83"
Set-TextValue "E3" "  -5.34%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "561.28"
Set-TextValue "E5" "  -3.75%  "
Set-TextValue "D6" "182.32"
Set-TextValue "E6" "  -6.64%  "
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "D8" "0.589"
Set-TextValue "E8" "  -3.11%  "
Set-TextValue "D9" "3.336.92"
Set-TextValue "E9" "  -5.15%  "
Set-TextValue "D10" "0.185"
Set-TextValue "E10" "  -8.78%  "
Set-TextValue "D11" "0.586"
Set-TextValue "E11" "  -6.76%  "
Set-TextValue "D12" "47.30"
Set-TextValue "E12" "  -8.35%  "
Set-TextValue "D13" "0.0000266"
Set-TextValue "E13" "  -6.92%  "
Set-TextValue "D14" "3.872.29"
Set-TextValue "E14" "  -5.62%  "
Set-TextValue "D15" "8.58"
Set-TextValue "E15" "  -6.64%  "
Set-TextValue "D16" "603.78"
Set-TextValue "E16" "  -9.26%  "
Set-TextValue "D17" "18.16"
Set-TextValue "E17" "  -1.44%  "
Set-TextValue "D18" "66.234.49"
Set-TextValue "E18" "  -4.61%  "
Set-TextValue "D19" "3.341.24"
Set-TextValue "E19" "  -5.81%  "
Set-TextValue "E20" "  -3.91%  "
Set-TextValue "D21" "11.41"
Set-TextValue "E21" "  -8.74%  "
Set-TextValue "D22" "0.905"
Set-TextValue "E22" "  -5.98%  "
Set-TextValue "D23" "16.83"
Set-TextValue "E23" "  -8.18%  "
Set-TextValue "D24" "5.07"
Set-TextValue "E24" "  -4.64%  "
Set-TextValue "D25" "100.48"
Set-TextValue "E25" "  -3.75%  "
Set-TextValue "D26" "4.03"
Set-TextValue "E26" "  -7.64%  "
Set-TextValue "E27" "  +0.18%  "
Set-TextValue "D28" "2.66"
Set-TextValue "E28" "  -7.94%  "
Set-TextValue "D29" "9.31"
Set-TextValue "E29" "  -8.09%  "
Set-TextValue "D30" "8.70"
Set-TextValue "E30" "  -9.20%  "
Set-TextValue "D31" "30.63"
Set-TextValue "E31" "  -7.36%  "
Set-TextValue "D32" "6.24"
Set-TextValue "E32" "  -7.39%  "
Set-TextValue "D33" "3.75"
Set-TextValue "E33" "  -14.45%  "
Set-TextValue "D34" "11.00"
Set-TextValue "E34" "  -6.20%  "
Set-TextValue "E35" "  -5.51%  "
Set-TextValue "D36" "3.787.28"
Set-TextValue "E36" "  -0.03%  "
Set-TextValue "B37" "Bittensor"
Set-TextValue "C37" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D37" "534.95"
Set-TextValue "E37" "  +6.93%  "
Set-TextValue "B38" "OKB"
Set-TextValue "C38" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D38" "57.92"
Set-TextValue "E38" "  -6.35%  "
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  +0.40%  "
Set-TextValue "D40" "3.47"
Set-TextValue "E40" "  -6.03%  "
Set-TextValue "D41" "0.0₃0712"
Set-TextValue "E41" "  -12.16%  "
Set-TextValue "D42" "2.66"
Set-TextValue "E42" "  -8.51%  "
Set-TextValue "D43" "0.125"
Set-TextValue "E43" "  -7.00%  "
Set-TextValue "D44" "0.340"
Set-TextValue "E44" "  -8.19%  "
Set-TextValue "D45" "31.84"
Set-TextValue "E45" "  -7.60%  "
Set-TextValue "D46" "0.0413"
Set-TextValue "E46" "  -7.75%  "
Set-TextValue "B47" "ApeXProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "3.24"
Set-TextValue "E47" "  -4.40%  "
Set-TextValue "B48" "CoreDAO"
Set-TextValue "C48" "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D48" "3.06"
Set-TextValue "E48" "  +15.25%  "
Set-TextValue "B49" "ThetaToken"
Set-TextValue "C49" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D49" "2.61"
Set-TextValue "E49" "  -8.58%  "
Set-TextValue "B50" "Stellar"
Set-TextValue "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.129"
Set-TextValue "E50" "  -4.98%  "
Set-TextValue "D51" "0.999"
Set-TextValue "E51" "  -0.27%  "